$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A8").Value = "     Développement et test de la classe Outils"
$ws.Range("A11").Value = "     Développement et test de la classe Reservation"
$ws.Range("A14").Value = "     Développement et test de la classe Utilisateur"
$ws.Range("A17").Value = "     Développement et test de la classe Salle"
$ws.Range("A20").Value = "     Développement et test méthode __construct"
$ws.Range("A22").Value = "     Développement et test méthode annulerReservation"
$ws.Range("A21").Value = "     Développement et test méthode __destruct"
$ws.Range("A23").Value = "     Développement et test méthode aPasseDesReservations"
$ws.Range("A24").Value = "     Développement et test méthode confirmerReservation"
$ws.Range("A25").Value = "     Développement et test méthode creerLesDigicodesManquants"
$ws.Range("A26").Value = "     Développement et test méthode creerUtilisateur"
$ws.Range("A27").Value = "     Développement et test méthode envoyerMdp"
$ws.Range("A28").Value = "     Développement et test méthode estLeCreateur"
$ws.Range("A30").Value = "     Développement et test méthode existeUtilisateur"
$ws.Range("A29").Value = "     Développement et test méthode existeReservation"
$ws.Range("A31").Value = "     Développement et test méthode genererUnDigicode"
$ws.Range("A32").Value = "     Développement et test méthode getLesReservations"
$ws.Range("A34").Value = "     Développement et test méthode getNiveauUtilisateur"
$ws.Range("A33").Value = "     Développement et test méthode getLesSalles"
$ws.Range("A35").Value = "     Développement et test méthode getReservation"
$ws.Range("A36").Value = "     Développement et test méthode getUtilisateur"
$ws.Range("A37").Value = "     Développement et test méthode modifierMdpUser"
$ws.Range("A38").Value = "     Développement et test méthode supprimerUtilisateur"
$ws.Range("A39").Value = "     Développement et test méthode testerDigicodeBatiment"
$ws.Range("A40").Value = "     Développement et test méthode testerDigicodeSalle"
$ws.Range("A45").Value = "     Développement et test page index"
$ws.Range("A47").Value = "     Développement et test CtrlAnnulerReservation"
$ws.Range("A48").Value = "     Développement et test VueAnnulerReservation"
$ws.Range("A50").Value = "     Développement et test CtrlChangerDeMdp"
$ws.Range("A51").Value = "     Développement et test VueChangerDeMdp"
$ws.Range("A53").Value = "     Développement et test CtrlConfirmerReservation"
$ws.Range("A54").Value = "     Développement et test VueConfirmerReservation"
$ws.Range("A56").Value = "     Développement et test CtrlConnecter"
$ws.Range("A57").Value = "     Développement et test VueConnecter"
$ws.Range("A59").Value = "     Développement et test CtrlConsulterReservations"
$ws.Range("A60").Value = "     Développement et test VueConsulterReservations"
$ws.Range("A62").Value = "     Développement et test CtrlConsulterSalles"
$ws.Range("A63").Value = "     Développement et test VueConsulterSalles"
$ws.Range("A65").Value = "     Développement et test CtrlCreerUtilisateur"
$ws.Range("A66").Value = "     Développement et test VueCreerUtilisateur"
$ws.Range("A68").Value = "     Développement et test CtrlDemanderMdp"
$ws.Range("A69").Value = "     Développement et test VueDemanderMdp"
$ws.Range("A71").Value = "     Développement et test CtrlMenu"
$ws.Range("A72").Value = "     Développement et test VueMenu"
$ws.Range("A74").Value = "     Développement et test CtrlSupprimerUtilisateur"
$ws.Range("A75").Value = "     Développement et test VueSupprimerUtilisateur"
$ws.Range("A77").Value = "     Développement et test CtrlTelechargerApk"

$ws.Range("A5").Select()